$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "magapoke_2025-12-10"
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

$ws.Cells.Item(1,1).Value = "rank"
$ws.Cells.Item(1,2).Value = "title"

$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "歪みの虜"
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "MYS"
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "黒月のイェルクナハト"
$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = "スルガメテオ"
$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = "ドリーム☆ジャンボ☆ガール"
$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = "ハードワーカー中田"
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "アイドラトリィ"
$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = "K-9~警視庁公安部公安第9課異能対策係~"
$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = "その青春"
$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = "黄昏町プリズナーズ"
$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = "せいぶつ部の田辺くん"
$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = "篝家の８兄弟"
$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = "ともだちづくり"
$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = "ナキナギ"
$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = "ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜"
$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = "追放されなかった男　～二度目の人生は土下座から始まりました～"
$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = "眠れる森のレガ"
$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = "春くらり"
$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = "普通の本はありません！"
$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(21,2).Value = "夜鐘のキト"
$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(22,2).Value = "お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！"
$ws.Cells.Item(23,1).Value = 22
$ws.Cells.Item(23,2).Value = "ゼロとヒャク"
$ws.Cells.Item(24,1).Value = 23
$ws.Cells.Item(24,2).Value = "それがメイドのカンナです"
$ws.Cells.Item(25,1).Value = 24
$ws.Cells.Item(25,2).Value = "屋根の下のアルテミス"
$ws.Cells.Item(26,1).Value = 25
$ws.Cells.Item(26,2).Value = "花子狩り"
$ws.Cells.Item(27,1).Value = 26
$ws.Cells.Item(27,2).Value = "皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～"
$ws.Cells.Item(28,1).Value = 27
$ws.Cells.Item(28,2).Value = "ハプスブルク家の華麗なる受難"
$ws.Cells.Item(29,1).Value = 28
$ws.Cells.Item(29,2).Value = "白銀のキュイジーヌ～明治外交官の料理人～"
$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = "生きたがりの人狼"
$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = "人生逆転ダンジョン"
$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = "限界集落を脱村した錬金術士、都会で""最強""なのがバレまくる。～老害どもにはいい加減愛想が尽きました～"
$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = "ハナバス　苔石花江のバスケ論"
$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = "永久のユウグレ"
$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = "君が監督！"
$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = "明智ナンバーワン"
$ws.Cells.Item(37,1).Value = 36
$ws.Cells.Item(37,2).Value = "平成転生"
$ws.Cells.Item(38,1).Value = 37
$ws.Cells.Item(38,2).Value = "鳴るさんだぁ"
$ws.Cells.Item(39,1).Value = 38
$ws.Cells.Item(39,2).Value = "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～"
$ws.Cells.Item(40,1).Value = 39
$ws.Cells.Item(40,2).Value = "じゅーくぼっくす"
$ws.Cells.Item(41,1).Value = 40
$ws.Cells.Item(41,2).Value = "〈小市民〉 春期限定いちごタルト事件"
$ws.Cells.Item(42,1).Value = 41
$ws.Cells.Item(42,2).Value = "JK Biker"
$ws.Cells.Item(43,1).Value = 42
$ws.Cells.Item(43,2).Value = "卒業アルバムの彼女たち"
$ws.Cells.Item(44,1).Value = 43
$ws.Cells.Item(44,2).Value = "ナマイキ旭ちゃんをわからせたい"
$ws.Cells.Item(45,1).Value = 44
$ws.Cells.Item(45,2).Value = "イエティ、とある日々"
$ws.Cells.Item(46,1).Value = 45
$ws.Cells.Item(46,2).Value = "東京デスレース"

$src = $ws1.Range("A1:B1")
$src.Copy()
$dst = $ws.Range("A1:B1")
$dst.PasteSpecial(-4122)
$ws.Range("A1").Select()
$ws1.Activate()
